# Apply corrections to PNAD 2009 "agressao" sheet:
# 1. Fix the header on B2 (was the stray pandas "unnamed: 1_level_1" label,
#    should read the same "total" label as B1).
# 2. Remove the two empty category-header rows ("situação do domicílio" and
#    "grandes regiões e unidades da federação") that had no data underneath
#    them, which shifts every following data row up by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the mis-labelled header cell.
$ws.Range("B2").Value = "total"

# 2) Delete the "grandes regiões e unidades da federação" row (row 8) first,
#    then the "situação do domicílio" row (row 5), deleting the lower one
#    first so the earlier row index stays valid.
$ws.Rows("8:8").Delete() | Out-Null
$ws.Rows("5:5").Delete() | Out-Null
